$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$companies = @(
    "Wash World",
    "Carve",
    "GroupM",
    "Brøndbyernes I.F.",
    "Meew",
    "Funelo",
    "PreCure",
    "Wilke",
    "OOONO",
    "Elbek & Vejrup",
    "Ellab",
    "Lejka",
    "Firi"
)

for ($i = 0; $i -lt $companies.Length; $i++) {
    $row = 22 + $i
    $ws.Cells.Item($row, 3).Value = $companies[$i]
}

Write-Output "done"
